$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns with latest scraped values.
# Force text number format so values are stored as text (matching original inlineStr cells)
# rather than being auto-converted to numeric/percentage values by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.47%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.91%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.736"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.63%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08018"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.30%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.510"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.95%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.618"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.63%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.938"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.58%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.63%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9175"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.40%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1249"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1946"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.21%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.717"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "16.62%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09255"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.35%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03561"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.76%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1050"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.52%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001294"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.88%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006291"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.18%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.368"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.06%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.36%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1376"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.73%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2672"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.32%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04451"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.97%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004424"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.28%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02527"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.28%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05434"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.37%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007545"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.39%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009920"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "12.35%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1401"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.05%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.76%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01125"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "16.72%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006796"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.94%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003043"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.62%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002283"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.87%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
